$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 108, shifting rows 108:162 down to 109:163
$ws.Rows("108:108").Insert()

# Fill the new row 108 with values (copy static fields from former row 108 / now row 109)
$ws.Range("A108").Value = 4
$ws.Range("B108").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C108").Value = "Los Lagos"
$ws.Range("D108").Value = 44466
$ws.Range("E108").Value = 10
$ws.Range("F108").Value = 100112037
$ws.Range("G108").Value = "Cebollín"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 90
$ws.Range("K108").Value = 6000
$ws.Range("L108").Value = 6000
$ws.Range("M108").Value = 6000
$ws.Range("N108").Value = "$/paquete 36 unidades"
$ws.Range("O108").Value = "Región Metropolitana"
$ws.Range("P108").Value = 167
$ws.Range("Q108").Value = 36
$ws.Range("R108").Value = "Hortaliza"
